# Refresh cryptocurrency price/volume data (GitHub Actions scheduled update).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell reference -> new value. Values that are purely numeric-looking
# are prefixed with a leading apostrophe so Excel stores them as literal text
# (preserving formatting such as trailing zeros), matching the source data which
# is stored as text rather than as numbers.
$updates = @(
    @{ Cell = "D2"; Value = "30.050.35" }
    @{ Cell = "E2"; Value = "  -1.39%  " }
    @{ Cell = "D3"; Value = "1.861.61" }
    @{ Cell = "E3"; Value = "  -2.66%  " }
    @{ Cell = "E4"; Value = "  +0.24%  " }
    @{ Cell = "D5"; Value = "'233.97" }
    @{ Cell = "D6"; Value = "'1.003" }
    @{ Cell = "E6"; Value = "  +0.35%  " }
    @{ Cell = "D7"; Value = "'0.4671" }
    @{ Cell = "E7"; Value = "  -2.26%  " }
    @{ Cell = "D8"; Value = "'0.2820" }
    @{ Cell = "E8"; Value = "  -0.92%  " }
    @{ Cell = "D9"; Value = "'0.06546" }
    @{ Cell = "E9"; Value = "  -2.42%  " }
    @{ Cell = "D10"; Value = "'20.20" }
    @{ Cell = "E10"; Value = "  +6.80%  " }
    @{ Cell = "D11"; Value = "'0.07761" }
    @{ Cell = "E11"; Value = "  +0.82%  " }
    @{ Cell = "D12"; Value = "'96.19" }
    @{ Cell = "E12"; Value = "  -6.12%  " }
    @{ Cell = "D13"; Value = "1.877.40" }
    @{ Cell = "E13"; Value = "  -1.84%  " }
    @{ Cell = "D14"; Value = "'5.054" }
    @{ Cell = "E14"; Value = "  -3.07%  " }
    @{ Cell = "D15"; Value = "'0.6677" }
    @{ Cell = "E15"; Value = "  -0.56%  " }
    @{ Cell = "D16"; Value = "'282.03" }
    @{ Cell = "E16"; Value = "  +3.59%  " }
    @{ Cell = "D17"; Value = "30.057.51" }
    @{ Cell = "E17"; Value = "  -1.49%  " }
    @{ Cell = "D18"; Value = "'1.002" }
    @{ Cell = "E18"; Value = "  +0.11%  " }
    @{ Cell = "D19"; Value = "2.119.38" }
    @{ Cell = "E19"; Value = "  -1.85%  " }
    @{ Cell = "D20"; Value = "'12.51" }
    @{ Cell = "E20"; Value = "  -1.45%  " }
    @{ Cell = "D21"; Value = "'5.339" }
    @{ Cell = "E21"; Value = "  -1.70%  " }
    @{ Cell = "B22"; Value = "BinanceUSD" }
    @{ Cell = "C22"; Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd" }
    @{ Cell = "D22"; Value = "'1.005" }
    @{ Cell = "E22"; Value = "  +0.47%  " }
    @{ Cell = "B23"; Value = "ShibaInu" }
    @{ Cell = "C23"; Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib" }
    @{ Cell = "D23"; Value = "'0.000007214" }
    @{ Cell = "E23"; Value = "  -3.46%  " }
    @{ Cell = "D24"; Value = "'6.131" }
    @{ Cell = "E24"; Value = "  -2.82%  " }
    @{ Cell = "B25"; Value = "Monero" }
    @{ Cell = "C25"; Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr" }
    @{ Cell = "D25"; Value = "'167.10" }
    @{ Cell = "E25"; Value = "  +0.19%  " }
    @{ Cell = "B26"; Value = "Cosmos" }
    @{ Cell = "C26"; Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom" }
    @{ Cell = "D26"; Value = "'9.289" }
    @{ Cell = "E26"; Value = "  -1.11%  " }
    @{ Cell = "D27"; Value = "'18.92" }
    @{ Cell = "E27"; Value = "  -2.31%  " }
    @{ Cell = "D28"; Value = "'1.949" }
    @{ Cell = "E28"; Value = "  -5.45%  " }
    @{ Cell = "D29"; Value = "'1.372" }
    @{ Cell = "E29"; Value = "  -0.67%  " }
    @{ Cell = "D30"; Value = "'0.09640" }
    @{ Cell = "E30"; Value = "  -4.09%  " }
    @{ Cell = "D31"; Value = "'4.368" }
    @{ Cell = "E31"; Value = "  -5.45%  " }
    @{ Cell = "D32"; Value = "'1.467" }
    @{ Cell = "E32"; Value = "  -2.97%  " }
    @{ Cell = "D33"; Value = "'4.071" }
    @{ Cell = "E33"; Value = "  -3.77%  " }
    @{ Cell = "D34"; Value = "'0.04645" }
    @{ Cell = "E34"; Value = "  -1.90%  " }
    @{ Cell = "D35"; Value = "'0.6977" }
    @{ Cell = "E35"; Value = "  -3.86%  " }
    @{ Cell = "D36"; Value = "'1.082" }
    @{ Cell = "E36"; Value = "  -2.52%  " }
    @{ Cell = "B37"; Value = "HuobiToken" }
    @{ Cell = "C37"; Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht" }
    @{ Cell = "D37"; Value = "'2.709" }
    @{ Cell = "E37"; Value = "  -0.28%  " }
    @{ Cell = "B38"; Value = "VeChain" }
    @{ Cell = "C38"; Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet" }
    @{ Cell = "D38"; Value = "'0.01852" }
    @{ Cell = "E38"; Value = "  -3.74%  " }
    @{ Cell = "B39"; Value = "FraxShare" }
    @{ Cell = "C39"; Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs" }
    @{ Cell = "D39"; Value = "'6.402" }
    @{ Cell = "E39"; Value = "  +2.02%  " }
    @{ Cell = "B40"; Value = "MXToken" }
    @{ Cell = "C40"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" }
    @{ Cell = "D40"; Value = "'2.515" }
    @{ Cell = "E40"; Value = "  -3.66%  " }
    @{ Cell = "B41"; Value = "Aave" }
    @{ Cell = "C41"; Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave" }
    @{ Cell = "D41"; Value = "'71.26" }
    @{ Cell = "E41"; Value = "  -4.89%  " }
    @{ Cell = "B42"; Value = "TrustWalletToken" }
    @{ Cell = "C42"; Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt" }
    @{ Cell = "D42"; Value = "'0.8590" }
    @{ Cell = "E42"; Value = "  +0.10%  " }
    @{ Cell = "B43"; Value = "RenderToken" }
    @{ Cell = "C43"; Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr" }
    @{ Cell = "D43"; Value = "'1.939" }
    @{ Cell = "E43"; Value = "  -1.63%  " }
    @{ Cell = "B44"; Value = "PaxDollar" }
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp" }
    @{ Cell = "D44"; Value = "'1.003" }
    @{ Cell = "E44"; Value = "  +0.44%  " }
    @{ Cell = "B45"; Value = "Quant" }
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt" }
    @{ Cell = "D45"; Value = "'102.84" }
    @{ Cell = "E45"; Value = "  -2.08%  " }
    @{ Cell = "B46"; Value = "TheSandbox" }
    @{ Cell = "C46"; Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand" }
    @{ Cell = "D46"; Value = "'0.4157" }
    @{ Cell = "E46"; Value = "  -2.53%  " }
    @{ Cell = "B47"; Value = "Maker" }
    @{ Cell = "C47"; Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr" }
    @{ Cell = "D47"; Value = "'976.18" }
    @{ Cell = "E47"; Value = "  +6.40%  " }
    @{ Cell = "B48"; Value = "Aptos" }
    @{ Cell = "C48"; Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt" }
    @{ Cell = "D48"; Value = "'7.171" }
    @{ Cell = "E48"; Value = "  -3.60%  " }
    @{ Cell = "B49"; Value = "EnergySwap" }
    @{ Cell = "C49"; Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens" }
    @{ Cell = "D49"; Value = "'9.062" }
    @{ Cell = "E49"; Value = "  +2.62%  " }
    @{ Cell = "B50"; Value = "Elrond" }
    @{ Cell = "C50"; Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld" }
    @{ Cell = "D50"; Value = "'33.78" }
    @{ Cell = "E50"; Value = "  -2.96%  " }
    @{ Cell = "B51"; Value = "Algorand" }
    @{ Cell = "C51"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo" }
    @{ Cell = "D51"; Value = "'0.1139" }
    @{ Cell = "E51"; Value = "  -5.24%  " }
)

foreach ($update in $updates) {
    $ws.Range($update.Cell).Value = $update.Value
}

